$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "TONCOIN/USDT" -> "TON/USDT" (shared string referenced by cell A6)
$ws.Range("A6").Value = "TON/USDT"

# The previously-applied (now pointless) fill style on column A is removed,
# which also means cells A1:A9 drop their explicit style index.
$ws.Columns(1).ClearFormats()

# Cells below the table that ClearFormats() touched but that hold no data
# must stay empty/absent, not turn into blank-but-styled cells.
$ws.Range("A10:A12").ClearContents()

# Active selection moves from C13 to A7.
$ws.Range("A7").Select() | Out-Null
